$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new
# header cells I1:J1 before writing their text, so the new columns match
# the look of the existing header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels for the two new columns.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data, row by row (row 1 is the header handled above).
$data = @{
    2  = @(2, 3)
    3  = @(6, 7)
    4  = @(5, 5)
    5  = @(2, 5)
    6  = @(8, 8)
    7  = @(8, 8)
    8  = @(8, 8)
    9  = @(7, 8)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(10, 11)
    15 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
